$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("best-feasible-slns")

# Update sheet view: remove topLeftCell, change selection to O21
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O21").Select()

# Data updates in column C for rows 152-215
$values = @{
    152 = 117075
    153 = 118521
    154 = 118270
    155 = 118063
    156 = 115355
    157 = 118572
    158 = 119139
    159 = 117480
    160 = 116904
    161 = 118028
    162 = 216414
    163 = 217799
    164 = 216809
    165 = 216375
    166 = 212714
    167 = 213816
    168 = 217089
    169 = 219148
    170 = 213129
    171 = 220081
    172 = 303626
    173 = 301811
    174 = 301895
    175 = 300423
    176 = 304011
    177 = 301042
    178 = 304514
    179 = 295195
    180 = 300773
    181 = 306622
    182 = 21772
    183 = 21492
    184 = 20675
    185 = 21411
    186 = 21733
    187 = 22176
    188 = 21699
    189 = 21169
    190 = 22386
    191 = 20944
    192 = 40575
    193 = 41199
    194 = 41473
    195 = 40972
    196 = 40872
    197 = 41058
    198 = 40887
    199 = 42719
    200 = 42230
    201 = 41524
    202 = 57494
    203 = 59997
    204 = 57977
    205 = 60776
    206 = 58816
    207 = 59830
    208 = 58132
    209 = 58944
    210 = 58856
    211 = 60574
    212 = 56572
    213 = 58064
    214 = 55722
    215 = 58451
}

foreach ($row in ($values.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
